$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

foreach ($r in 3..14) {
    # Column H: decrement the numeric "period to expire" value by 1
    $hCell = $ws.Cells.Item($r, 8)
    $hCell.Value = $hCell.Value2 - 1

    # Column I: update the "last update" date text from 03-Nov-2025 to 04-Nov-2025
    # (kept as literal text, not an Excel date serial, matching the original cell type)
    $iCell = $ws.Cells.Item($r, 9)
    $iCell.Formula = '=TEXT("04-Nov-2025","@")'
    $iCell.Copy()
    $iCell.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = 0
